$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9999717487778564
$ws.Range("E2").Value = 0.9999717487778564

# Row 3
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = [double]"6.998151979627356E-07"
$ws.Range("E3").Value = [double]"6.998151979627356E-07"

# Row 4
$ws.Range("D4").Value = 0.9999991809706881
$ws.Range("E4").Value = 0.9999991809706881

# Row 5
$ws.Range("D5").Value = 0.05419778049529315
$ws.Range("E5").Value = 0.05419778049529315

# Row 6
$ws.Range("D6").Value = [double]"6.80692856056736E-16"
$ws.Range("E6").Value = [double]"6.80692856056736E-16"

# Row 8
$ws.Range("D8").Value = [double]"1.706659356146565E-05"
$ws.Range("E8").Value = 0.9999829334064385

# Row 9
$ws.Range("D9").Value = [double]"7.046085451998499E-239"

# Row 10
$ws.Range("D10").Value = [double]"2.048792234315811E-06"
$ws.Range("E10").Value = 0.9999979512077657

# Row 11
$ws.Range("D11").Value = 0.9999358588185379
$ws.Range("E11").Value = [double]"6.414118146214065E-05"
$ws.Range("F11").Value = 59.69873046875
$ws.Range("G11").Value = 0.5
